$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 22 (Squilla mantis, 2-RAP): weight and number corrected
$ws.Range("G22").Value = 0.253
$ws.Range("H22").Value = 11

# Remove the duplicate "Squilla mantis" row (row 47) entirely, shifting
# the remaining rows (Stones NA ... Wood NA) up by one.
$ws.Rows(47).Delete()
